$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing bug rows with their new descriptions.
$ws.Range("C23").Value = "need to make the xls file ignore populations if they are blank."
$ws.Range("C25").Value = "Write Vignette with nested source() calls, and create code_tree figure"
$ws.Range("C26").Value = "Testing of sumby"

# Remove the two trailing rows that are no longer needed.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()

# Match the saved selection state from the target workbook.
$ws.Range("C29").Select()
